$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Make room for one new line of explanation text under the "Example #2"
# formula-safety notes (pushes "Example #1" and everything below it down by
# one row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Insert()

# Row 13: lightly reworded bullet ("under TBS block" -> "under a TBS block").
$ws.Range("B13").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."

# Row 14 (brand new line, inherits style from the row above via the insert): continuation sentence.
$ws.Range("B14").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."

# Row 15: the old "reference" bullet slides down into this row, unchanged text.
$ws.Range("B15").Value = "* If a formula uses a reference to a cell that has moved during the merge, then the reference will not be arraged to be the new cell reference. "

# Row 16 (new row created by the insert above): the old "picture" bullet slides down into this row, unchanged text.
$ws.Range("B16").Value = "* You cannot change picture using ""ope=changepic"". This is because drawing information are not saved directly in the sheet."

# ---------------------------------------------------------------------------
# New "Score" column (E) + a "Total:" row above the example#1 table.
# ---------------------------------------------------------------------------

# "Total:" label, right aligned, plain style.
$ws.Range("D19").Value = "Total:"
$ws.Range("D19").HorizontalAlignment = -4152

# Bold SUM formula cell, one-decimal number format.
$ws.Range("E19").Formula = "=SUM(E21:E2000)"
$ws.Range("E19").NumberFormat = "#,##0.0"
$ws.Range("E19").Font.Bold = $true

# "Score" header cell - reuse the same look as the other header cells (B20:D20).
$ws.Range("B20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Score"

# TBS placeholder cell for the score value - reuse the bordered look of the
# other placeholder cells (B21:D21), then apply the number format/alignment.
$ws.Range("C21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152

# Match the saved selection state.
$ws.Range("E20").Select()

Write-Host "done"
